# Branch wise stock status added
# Corrects the Item Name / UOM pairings for several rows in the NoStock
# sheet (Ketonic / Kynol / Dinafex / Zithrox groups were mismatched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Dinafex 120mg / 60mg rows were swapped
$ws.Range("D4").Value = "Dinafex 60mg Tablet"
$ws.Range("D5").Value = "Dinafex 120mg Tablet"

# Ketonic rows: row 14 and row 16 item/UOM pairs were swapped
$ws.Range("D14").Value = "Ketonic 10mg Tablet"
$ws.Range("E14").Value = "20's"
$ws.Range("D16").Value = "Ketonic 30mg Injection"
$ws.Range("E16").Value = "5 's"

# Kynol rows 17-19: item/UOM pairs rotated
$ws.Range("D17").Value = "Kynol TR 100mg Capsule"
$ws.Range("E17").Value = "50 's"
$ws.Range("D18").Value = "Kynol TR 200mg Capsule"
$ws.Range("E18").Value = "30 's"
$ws.Range("D19").Value = "Kynol D 25mg Tablet"
$ws.Range("E19").Value = "60 's"

# Zithrox rows 24-27: item/UOM pairs rotated
$ws.Range("D24").Value = "Zithrox 15ml Suspension"
$ws.Range("E24").Value = "15 ml"
$ws.Range("D25").Value = "Zithrox 250mg Tablet - 6's"
$ws.Range("E25").Value = "6's"
$ws.Range("D26").Value = "Zithrox 500mg Tablet"
$ws.Range("E26").Value = "6 's"
$ws.Range("D27").Value = "Zithrox 30ml Dry Suspension"
$ws.Range("E27").Value = "30ml"
